$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing cell corrections (Q4 third GDP revision)
$ws.Range("H98").Value = 904864
$ws.Range("M99").Value = 8586173
$ws.Range("M100").Value = 8552546
$ws.Range("M102").Value = 8469972
$ws.Range("M104").Value = 8353736
$ws.Range("L105").Value = 82.48
$ws.Range("M105").Value = 8358613
$ws.Range("M109").Value = 8212234
$ws.Range("H132").Value = 935885
$ws.Range("M132").Value = 14067432
$ws.Range("L133").Value = 99.53
$ws.Range("M133").Value = 14469885
$ws.Range("M134").Value = 14605520
$ws.Range("L135").Value = 100.72
$ws.Range("M135").Value = 14690957
$ws.Range("M136").Value = 14867416
$ws.Range("M137").Value = 14861934
$ws.Range("M138").Value = 14671326
$ws.Range("M139").Value = 14467551
$ws.Range("D140").Value = 823650
$ws.Range("M140").Value = 14132331
$ws.Range("N140").Value = 8935486
$ws.Range("M141").Value = 13894627
$ws.Range("N141").Value = 8782001
$ws.Range("M142").Value = 13750534
$ws.Range("N142").Value = 8776269
$ws.Range("M143").Value = 13622403
$ws.Range("N143").Value = 8836711
$ws.Range("N144").Value = 9042519
$ws.Range("N145").Value = 9395139
$ws.Range("N146").Value = 9716003
$ws.Range("N147").Value = 9916834
$ws.Range("N148").Value = 10176467
$ws.Range("N149").Value = 10402048
$ws.Range("N150").Value = 10586448
$ws.Range("N151").Value = 10834361
$ws.Range("D186").Value = 590469
$ws.Range("H186").Value = 1208765
$ws.Range("M186").Value = 13889555
$ws.Range("N186").Value = 8211274
$ws.Range("M187").Value = 13965788
$ws.Range("N187").Value = 8231265
$ws.Range("M188").Value = 14177021
$ws.Range("N188").Value = 8311051
$ws.Range("M189").Value = 14260329
$ws.Range("N189").Value = 8315411
$ws.Range("L190").Value = 126.78
$ws.Range("M190").Value = 14422656
$ws.Range("N190").Value = 8370992
$ws.Range("D191").Value = 573526
$ws.Range("E191").Value = 207680
$ws.Range("H191").Value = 1073025
$ws.Range("M191").Value = 14527779
$ws.Range("N191").Value = 8368492
$ws.Range("O191").Value = 2551096
$ws.Range("M192").Value = 14532554
$ws.Range("N192").Value = 8290066
$ws.Range("O192").Value = 2547718
$ws.Range("M193").Value = 14761077
$ws.Range("N193").Value = 8365158
$ws.Range("O193").Value = 2572366
$ws.Range("M194").Value = 14843758
$ws.Range("N194").Value = 8388045
$ws.Range("O194").Value = 2570166
$ws.Range("M195").Value = 14888688
$ws.Range("N195").Value = 8325517
$ws.Range("O195").Value = 2571055
$ws.Range("L196").Value = 130.19
$ws.Range("M196").Value = 14967370
$ws.Range("N196").Value = 8320013
$ws.Range("O196").Value = 2563129
$ws.Range("M197").Value = 15089514
$ws.Range("N197").Value = 8341477
$ws.Range("O197").Value = 2566584
$ws.Range("M198").Value = 15141474
$ws.Range("N198").Value = 8330286
$ws.Range("O198").Value = 2567385
$ws.Range("M199").Value = 15309018
$ws.Range("N199").Value = 8356534
$ws.Range("O199").Value = 2587940
$ws.Range("M200").Value = 15440628
$ws.Range("N200").Value = 8350553
$ws.Range("O200").Value = 2599993
$ws.Range("E201").Value = 210444
$ws.Range("L201").Value = 133.48
$ws.Range("M201").Value = 15510402
$ws.Range("N201").Value = 8368296
$ws.Range("O201").Value = 2605571
$ws.Range("M202").Value = 15673712
$ws.Range("N202").Value = 8381853
$ws.Range("O202").Value = 2635483
$ws.Range("O203").Value = 2650976
$ws.Range("O204").Value = 2658798
$ws.Range("O205").Value = 2687347
$ws.Range("O206").Value = 2682401
$ws.Range("O207").Value = 2694752
$ws.Range("O208").Value = 2719335
$ws.Range("O209").Value = 2707777
$ws.Range("O210").Value = 2680342
$ws.Range("O211").Value = 2671979
$ws.Range("O212").Value = 2637237
$ws.Range("B604").Value = 11768425
$ws.Range("E604").Value = 164890
$ws.Range("F604").Value = 14906507
$ws.Range("G604").Value = 10615345
$ws.Range("H604").Value = 3891899
$ws.Range("M604").Value = 28616688
$ws.Range("O604").Value = 1773952
$ws.Range("B605").Value = 17631547
$ws.Range("D605").Value = 12201281
$ws.Range("E605").Value = 206118
$ws.Range("F605").Value = 66122705
$ws.Range("G605").Value = 55715631
$ws.Range("H605").Value = 18338675
$ws.Range("M605").Value = 44516865
$ws.Range("N605").Value = 18275696
$ws.Range("O605").Value = 1812672
$ws.Range("B606").Value = 9282648
$ws.Range("D606").Value = 5760519
$ws.Range("E606").Value = 263850
$ws.Range("F606").Value = 87028159
$ws.Range("G606").Value = 76904875
$ws.Range("H606").Value = 23754361
$ws.Range("M606").Value = 66265415
$ws.Range("N606").Value = 23685035
$ws.Range("O606").Value = 1925873
$ws.Range("B607").Value = 6650010
$ws.Range("D607").Value = 3163671
$ws.Range("E607").Value = 605297
$ws.Range("F607").Value = 77520547
$ws.Range("G607").Value = 75658920
$ws.Range("H607").Value = 22138430
$ws.Range("M607").Value = 86487995
$ws.Range("N607").Value = 26483306
$ws.Range("O607").Value = 2392485
$ws.Range("B608").Value = 5952294
$ws.Range("D608").Value = 2163842
$ws.Range("E608").Value = 950127
$ws.Range("F608").Value = 68598900
$ws.Range("G608").Value = 64868374
$ws.Range("H608").Value = 18662874
$ws.Range("M608").Value = 102836957
$ws.Range("N608").Value = 28174621
$ws.Range("O608").Value = 3184564
$ws.Range("B609").Value = 3754231
$ws.Range("E609").Value = 1066623
$ws.Range("F609").Value = 63217369
$ws.Range("G609").Value = 56553372
$ws.Range("H609").Value = 16358566
$ws.Range("M609").Value = 117204487
$ws.Range("N609").Value = 29141313
$ws.Range("O609").Value = 4115244
$ws.Range("B610").Value = 3499403
$ws.Range("D610").Value = 913794
$ws.Range("E610").Value = 3056261
$ws.Range("F610").Value = 49205690
$ws.Range("G610").Value = 46146982
$ws.Range("H610").Value = 13636582
$ws.Range("J610").Value = 109.41
$ws.Range("M610").Value = 128897007
$ws.Range("N610").Value = 29743917
$ws.Range("O610").Value = 7032268
$ws.Range("B611").Value = 3452808
$ws.Range("D611").Value = 762363
$ws.Range("E611").Value = 2435370
$ws.Range("F611").Value = 32121929
$ws.Range("G611").Value = 28459952
$ws.Range("H611").Value = 8377586
$ws.Range("M611").Value = 135364602
$ws.Range("N611").Value = 30141908
$ws.Range("O611").Value = 9326018
$ws.Range("B612").Value = 3435891
$ws.Range("E612").Value = 1329838
$ws.Range("F612").Value = 26459098
$ws.Range("G612").Value = 21945623
$ws.Range("H612").Value = 6701817
$ws.Range("M612").Value = 140240353
$ws.Range("N612").Value = 30517979
$ws.Range("O612").Value = 10527953
$ws.Range("B613").Value = 4225390
$ws.Range("D613").Value = 874958
$ws.Range("E613").Value = 1042516
$ws.Range("F613").Value = 23049165
$ws.Range("G613").Value = 20380357
$ws.Range("H613").Value = 6405718
$ws.Range("M613").Value = 143966724
$ws.Range("N613").Value = 30831500
$ws.Range("O613").Value = 11411558
$ws.Range("B614").Value = 4612118
$ws.Range("C614").Value = 5366.6
$ws.Range("D614").Value = 923321
$ws.Range("E614").Value = 834217
$ws.Range("F614").Value = 22888932
$ws.Range("G614").Value = 17873610
$ws.Range("H614").Value = 5746176
$ws.Range("I614").Value = 346.46
$ws.Range("J614").Value = 42.92
$ws.Range("L614").Value = 318.87
$ws.Range("M614").Value = 146765730
$ws.Range("N614").Value = 31018328
$ws.Range("O614").Value = 12093755

# New row 615 (new month of data appended)
$ws.Range("A615").Value = 44255
$ws.Range("B615").Value = 3457065
$ws.Range("C615").Value = 5719.333333333333
$ws.Range("D615").Value = 777465
$ws.Range("E615").Value = 1068857
$ws.Range("F615").Value = 19067992
$ws.Range("G615").Value = 16093068
$ws.Range("H615").Value = 5161117
$ws.Range("I615").Value = 349.35
$ws.Range("J615").Value = 43.23
$ws.Range("K615").Value = 15.69
$ws.Range("L615").Value = 318.76
$ws.Range("M615").Value = 149173800
$ws.Range("N615").Value = 31316363
$ws.Range("O615").Value = 12598111
